# Update the "想去人数" (want-to-go count) figures on both the
# "展览" and "全部类型" sheets, incrementing each by 1:
#   F2: 134 -> 135
#   F3: 1686 -> 1687
#   F6: 461 -> 462
#   F8: 75 -> 76

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 135
    $ws.Range("F3").Value = 1687
    $ws.Range("F6").Value = 462
    $ws.Range("F8").Value = 76
}
